# Bug fix for state saving
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (UID, Scholarship, Additional Feedback)
$data = @(
    @(1582872, "Test One", ""),
    @(1582872, "Cool Kids Club", ""),
    @(1536237, "Cool Kids Club", ""),
    @(1500178, "Cool Kids Club", ""),
    @(1500643, "Cool Kids Club", ""),
    @(1500178, "Test One", ""),
    @(1500643, "Test One", ""),
    @(1577306, "Test One", "")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
